# Generate Report for Handoff
# Refresh the localization-status report: the previous handoff
# (3abb78c5.../54f14ea3...) has been handed back and a brand new
# handoff round has just been generated (9be021f7.../ffff62b63081...).

$wb = $excel.ActiveWorkbook

$oldGuid1 = "3abb78c5-e03d-42e6-a5a1-4531fd5de059"
$oldGuid2 = "54f14ea3-ff94-4d72-8bbf-ea4d62077458"
$newGuid1 = "9be021f7-3b46-4926-baa6-f0185a86aa2f"
$newGuid2 = "ffff62b63081-f557-4714-a891-fd3339ead2c0"

$srcRepoBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/8191ff7835398d93ec5a2512e0a5ca224db7b29c/e2e/"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-13 05:17:19"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid1.md"
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("G2").Value = $newHoDate

$ov.Range("A3").Value = "$newGuid2.md"
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Range("G3").Value = $newHoDate

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "$srcRepoBase$newGuid1.md", "", "", "e2e\$newGuid1.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "$srcRepoBase$newGuid2.md", "", "", "e2e\$newGuid2.md")

$ov.Columns.Item(5).AutoFit()
$ov.Columns.Item(6).AutoFit()

# ---------------------------------------------------------------
# Language sheets ("zh-cn" / "de-de") share the same shape - loop.
# ---------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn"; XlfDate = "2016-08-13 05:17:12" },
    @{ Name = "de-de"; XlfDate = $newHoDate }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $ext = $lang.Name
    $newXlf = "$newGuid1.a72acc761a2beaeb55a0dbea5519726aba1b93b1.$ext.xlf"

    # Row 2 (was 3abb78c5...)
    $ws.Range("A2").Value = "$newGuid1.md"
    $ws.Range("C2").Value = $newStatus
    # Leading apostrophe forces literal text "True" instead of Excel
    # auto-coercing the recognised word into a boolean TRUE.
    $ws.Range("F2").Value = "'True"
    $ws.Range("G2").Value = $newXlf
    $ws.Range("H2").Value = $lang.XlfDate
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Row 3 (was 54f14ea3...)
    $ws.Range("A3").Value = "$newGuid2.md"
    $ws.Range("C3").Value = $newStatus
    $ws.Range("F3").Value = "'True"
    $ws.Range("G3").Value = $newXlf
    $ws.Range("H3").Value = $lang.XlfDate
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = "0001-01-01 00:00:00"

    # Clear the hyperlink look-and-feel left behind on I2/I3 now that
    # those cells no longer carry a "Latest Target File" hyperlink.
    $ws.Range("I2").Style = "Normal"
    $ws.Range("I3").Style = "Normal"

    # Hyperlinks: only column A (Source File Name) keeps a link; the
    # old column I (Latest Target File) links are dropped entirely.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$srcRepoBase$newGuid1.md", "", "", "$newGuid1.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "$srcRepoBase$newGuid2.md", "", "", "$newGuid2.md")

    $ws.Columns.Item(3).AutoFit()
    $ws.Columns.Item(9).AutoFit()
    $ws.Columns.Item(10).AutoFit()
}
